$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "About Me" paragraph had the word "switching" typed as three
# separate runs ("switch" | "ing" | " careers ... skill set.").
# A Find/Replace across that span collapses it back into a single run
# (identical text, identical formatting) - exactly what the target diff
# shows: the three runs become one run holding the whole sentence.
# -----------------------------------------------------------------------
$found1 = $d.Content.Find.Execute("switching", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "switching", 2)
if (-not $found1) {
    Write-Host "WARNING: 'switching' span not found for merge"
}

# -----------------------------------------------------------------------
# Change 2: the Portfolio link text
#   https://w-j-knight.github.io/william-knight.github.io/
# becomes
#   https://w-j-knight.github.io/william-knight  (existing run, shortened)
#   updateremuse                                  (new run)
#   .github.io/                                   (new run)
# i.e. the word "updateremuse" is inserted right after "william-knight"
# and before the trailing ".github.io/". The two new runs keep exactly
# the same visual formatting (InternetLink style, italic, same font/color)
# as the run they were split out of, so we briefly toggle a property to
# force the engine to keep them as distinct <w:r> elements instead of
# silently re-merging them with their identically-formatted neighbour.
# -----------------------------------------------------------------------
$full = $d.Content.Text
$marker = "william-knight.github.io"
$idx = $full.IndexOf($marker)
if ($idx -lt 0) {
    Write-Host "WARNING: portfolio URL marker not found"
} else {
    $insertPos = $idx + "william-knight".Length
    $insertion = "updateremuse"

    $insertionPoint = $d.Range($insertPos, $insertPos)
    $insertionPoint.InsertAfter($insertion)

    $newRunRange = $d.Range($insertPos, $insertPos + $insertion.Length)
    # Force a run boundary around the newly inserted text so it stays as
    # its own run(s) rather than being coalesced into the neighbouring
    # run that happens to share identical formatting.
    $newRunRange.Font.Bold = 1
    $newRunRange.Font.Bold = 0
}
